$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared-string rich text runs) ---
$ws.Range("A8").Value = 'Volume 31   Number  34'
$ws.Range("C9").Value = 'Report Covering the Week  8/19/2024  Through  8/25/2024'

# --- Numeric value updates across the precinct crime-stat table ---
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = -55.172413793103
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 1.666666666666
$ws.Range("L16").Value = 4.273504273504
$ws.Range("M16").Value = -46.491228070175
$ws.Range("N16").Value = -85.545023696682
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -30.434782608695
$ws.Range("I17").Value = 193
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = -5.392156862745
$ws.Range("L17").Value = 4.324324324324
$ws.Range("M17").Value = 42.962962962963
$ws.Range("N17").Value = -47.411444141689
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 25
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 218
$ws.Range("J18").Value = 173
$ws.Range("K18").Value = 26.011560693641
$ws.Range("L18").Value = 9
$ws.Range("M18").Value = -23.239436619718
$ws.Range("N18").Value = -75.255391600454
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 6.25
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = -30.985915492957
$ws.Range("I19").Value = 467
$ws.Range("J19").Value = 504
$ws.Range("K19").Value = -7.341269841269
$ws.Range("L19").Value = 7.603686635944
$ws.Range("M19").Value = 51.623376623376
$ws.Range("N19").Value = 36.151603498542
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 111
$ws.Range("K20").Value = -34.234234234234
$ws.Range("L20").Value = -33.027522935779
$ws.Range("M20").Value = -24.742268041237
$ws.Range("N20").Value = -87.102473498233
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -34.285714285714
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 165
$ws.Range("H21").Value = -32.121212121212
$ws.Range("I21").Value = 1089
$ws.Range("J21").Value = 1124
$ws.Range("K21").Value = -3.113879003558
$ws.Range("L21").Value = 3.027436140018
$ws.Range("M21").Value = 3.027436140018
$ws.Range("N21").Value = -64.212947748932
$ws.Range("L22").Value = -53.125
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 109
$ws.Range("J23").Value = 136
$ws.Range("K23").Value = -19.85294117647
$ws.Range("L23").Value = 5.825242718446
$ws.Range("M23").Value = 19.780219780219
$ws.Range("C24").Value = 29
$ws.Range("E24").Value = 45
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 3.296703296703
$ws.Range("I24").Value = 732
$ws.Range("J24").Value = 698
$ws.Range("K24").Value = 4.871060171919
$ws.Range("L24").Value = -9.629629629629
$ws.Range("M24").Value = -8.385481852315
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 11.764705882352
$ws.Range("I25").Value = 235
$ws.Range("J25").Value = 93
$ws.Range("K25").Value = 152.688172043011
$ws.Range("L25").Value = 32.022471910112
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 53
$ws.Range("H26").Value = -54.716981132075
$ws.Range("I26").Value = 343
$ws.Range("J26").Value = 329
$ws.Range("K26").Value = 4.255319148936
$ws.Range("L26").Value = 6.853582554517
$ws.Range("M26").Value = 8.888888888888
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("I28").Value = 44
$ws.Range("K28").Value = 10
$ws.Range("L28").Value = 46.666666666666
$ws.Range("N29").Value = -87.272727272727
$ws.Range("N30").Value = -87.755102040816
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 19
$ws.Range("K31").Value = 90
$ws.Range("L31").Value = 35.714285714285
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = 100
$ws.Range("L33").Value = 100

# --- Cells that flip from a number to the "0"/"***.*" text placeholder ---
# Use a quote-prefixed value to force text storage, then copy number
# formatting from an existing text-styled cell so the style index matches.
$zeroDonor = $ws.Range("D14")   # existing style-14 text cell holding "0"
$ws.Range("C20").Value = "'0"
$ws.Range("C22").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("D31").Value = "'0"
$zeroDonor.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)

$starDonor = $ws.Range("D14")   # same style-14 text cell, for the "***.*" cells
$ws.Range("E27").Value = "***.*"
$ws.Range("E28").Value = "***.*"
$ws.Range("E31").Value = "***.*"
$starDonor.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)

# --- Cells that flip from the "0" text placeholder back to a real number ---
$numDonor = $ws.Range("I33")    # existing style-15 numeric cell
$ws.Range("C33").Value = 1
$ws.Range("F33").Value = 1
$numDonor.Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("F33").PasteSpecial(-4122)

$excel.CutCopyMode = 0
